$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.947.63"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.155.48"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.37"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.50"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.155.60"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("E10").Value = "  -2.25%  "
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.501"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000264"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.29"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.674.00"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.915.70"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.161.09"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.15"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "505.49"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.91"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.716"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.21"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("E24").Value = "  -1.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.55"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.99"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.60%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.79"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.65"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.16%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.25"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.49"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.45%  "
$ws.Range("E36").Value = "  -1.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "484.85"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0891"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0417"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.93"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.77"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.997.48"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.116"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.02%  "
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("E45").Value = "  -4.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.25"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0588"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("D48").ClearFormats()
$ws.Range("E49").Value = "  -1.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.24"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.48"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +14.18%  "
